$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (style) of the last existing data row (328) down through
# the new rows that are being appended (329-343), without touching values.
$ws.Range("A328:D328").Copy() | Out-Null
$ws.Range("A329:D343").PasteSpecial(-4122) | Out-Null

$dates = @(44403, 44404, 44405, 44406, 44407, 44408, 44409, 44410, 44411, 44412, 44413, 44414, 44415, 44416, 44417)
$newPos = @(1, 1, 0, 1, 2, 0, 0, 6, 0, 0, 1, 0, 1, 2, 1)
$rolling = @(5, 6, 6, 6, 8, 6, 5, 10, 9, 9, 9, 7, 8, 10, 5)
$rollingPer100k = @(57.49770009199631, 68.99724011039559, 68.99724011039559, 68.99724011039559, 91.99632014719411, 68.99724011039559, 57.49770009199631, 114.9954001839926, 103.4958601655934, 103.4958601655934, 103.4958601655934, 80.49678012879485, 91.99632014719411, 114.9954001839926, 57.49770009199631)

$startRow = 329
for ($i = 0; $i -lt $dates.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $newPos[$i]
    $ws.Cells.Item($row, 3).Value = $rolling[$i]
    $ws.Cells.Item($row, 4).Value = $rollingPer100k[$i]
}
